$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.760.48'
$ws.Range("E2").Value = '  +3.00%  '
$ws.Range("D3").Value = '3.360.37'
$ws.Range("E3").Value = '  +8.82%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''262.62'
$ws.Range("E5").Value = '  +11.07%  '
$ws.Range("D6").Value = '''635.45'
$ws.Range("E6").Value = '  +5.42%  '
$ws.Range("D7").Value = '''1.38'
$ws.Range("E7").Value = '  +25.45%  '
$ws.Range("D8").Value = '''0.393'
$ws.Range("E8").Value = '  +4.00%  '
$ws.Range("D9").Value = '''0.999'
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '''0.869'
$ws.Range("E10").Value = '  +11.76%  '
$ws.Range("D11").Value = '3.355.09'
$ws.Range("E11").Value = '  +8.84%  '
$ws.Range("E12").Value = '  +2.42%  '
$ws.Range("D13").Value = '98.579.68'
$ws.Range("E13").Value = '  +3.55%  '
$ws.Range("D14").Value = '''36.27'
$ws.Range("E14").Value = '  +9.29%  '
$ws.Range("D15").Value = '''0.0000249'
$ws.Range("E15").Value = '  +5.66%  '
$ws.Range("D16").Value = '3.965.53'
$ws.Range("E16").Value = '  +8.40%  '
$ws.Range("D17").Value = '''5.55'
$ws.Range("E17").Value = '  +5.20%  '
$ws.Range("D18").Value = '3.362.57'
$ws.Range("E18").Value = '  +9.39%  '
$ws.Range("D19").Value = '''3.58'
$ws.Range("E19").Value = '  +2.76%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''15.13'
$ws.Range("E20").Value = '  +6.82%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '''494.18'
$ws.Range("E21").Value = '  +9.63%  '
$ws.Range("E22").Value = '  +10.69%  '
$ws.Range("D23").Value = '''0.0000212'
$ws.Range("E23").Value = '  +12.88%  '
$ws.Range("E24").Value = '  +8.28%  '
$ws.Range("E25").Value = '  +4.92%  '
$ws.Range("D26").Value = '''88.82'
$ws.Range("E26").Value = '  +4.93%  '
$ws.Range("D27").Value = '''12.08'
$ws.Range("E27").Value = '  +5.47%  '
$ws.Range("E28").Value = '  +8.81%  '
$ws.Range("D29").Value = '''0.284'
$ws.Range("E29").Value = '  +19.24%  '
$ws.Range("D30").Value = '''0.197'
$ws.Range("E30").Value = '  +11.46%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("E33").Value = '  +21.12%  '
$ws.Range("D34").Value = '''9.60'
$ws.Range("E34").Value = '  +9.09%  '
$ws.Range("D35").Value = '''27.81'
$ws.Range("E35").Value = '  +9.72%  '
$ws.Range("D36").Value = '''7.40'
$ws.Range("E36").Value = '  +2.34%  '
$ws.Range("E37").Value = '  +8.77%  '
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("D39").Value = '''508.29'
$ws.Range("E39").Value = '  +6.45%  '
$ws.Range("D40").Value = '''0.468'
$ws.Range("E40").Value = '  +9.70%  '
$ws.Range("D41").Value = '''24.84'
$ws.Range("E41").Value = '  +2.94%  '
$ws.Range("D42").Value = '''3.84'
$ws.Range("E42").Value = '  +5.68%  '
$ws.Range("D43").Value = '''1.28'
$ws.Range("E43").Value = '  +5.50%  '
$ws.Range("D44").Value = '''3.29'
$ws.Range("E44").Value = '  +6.64%  '
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("E46").Value = '  +15.00%  '
$ws.Range("D47").Value = '''161.76'
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("E48").Value = '  +6.42%  '
$ws.Range("D49").Value = '''46.52'
$ws.Range("E49").Value = '  +6.44%  '
$ws.Range("D50").Value = '''4.65'
$ws.Range("E50").Value = '  +10.01%  '
$ws.Range("D51").Value = '''1.37'
$ws.Range("E51").Value = '  +9.59%  '
